# Timetable Jonny - add working time rows + reformat the date column
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content, entered in the same order the original author typed it in
# (this keeps the shared-string table ordering identical to the real edit) ---

# Row 4 (new row): Duration + Task first
$ws.Range("B4").Value = "12:00-15:00"
$ws.Range("C4").Value = "Learning Git Hub"

# Row 2: re-type the date as plain text in a new format (not a real date)
$ws.Range("A2").Value = "31.04.2020"

# Row 5 (new row): Duration
$ws.Range("B5").Value = "11:00-16:00"

# Row 6 (new row): Duration
$ws.Range("B6").Value = "17:00-19:00"

# Row 5 (new row): Task
$ws.Range("C5").Value = "Setting up Group Paper, adding User Requirements, Project Goals, Workpackages, Naming conventions"

# Row 6 (new row): Task
$ws.Range("C6").Value = "Finishing first version of Group Paper, Setting up first Version of Wiki"

# --- Dates column: real date values with a left-aligned short-date format ---
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("A3").NumberFormat = "mm-dd-yy"
$ws.Range("A3").Value = 43922

# Copy A3's format onto the rest of the date column so they all share one style
$ws.Range("A3").Copy()
$ws.Range("A4:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Value = 43923
$ws.Range("A5").Value = 43925
$ws.Range("A6").Value = 43926

# --- Selection, matching where the author's cursor ended up ---
$ws.Range("C9").Select() | Out-Null
